# Fix for close price scraping from NSE: refresh LTP (col B) and PREV (col C)
# values on the "ltp" sheet for rows 2-26, restore the numeric style on B24
# (it now carries the same "#,##0.00" style as the other highlighted cells),
# and update the active selection left behind by the editor.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ltp")

# row -> (LTP, PREV)
$data = @{
    2  = @(621.9,     653)
    3  = @(3049.75,   3037.55)
    4  = @(474.1,     474.15)
    5  = @(1569.05,   1576.45)
    6  = @(6680,      6675.7)
    7  = @(189.9,     188.55)
    8  = @(263.9,     260.85)
    9  = @(47908.55,  48091.65)
    10 = @(825.65,    839.05)
    11 = @(4649.35,   4730.55)
    12 = @(162.75,    163.9)
    13 = @(1332.9,    1321.05)
    14 = @(653.7,     646.65)
    15 = @(1417.1,    1421)
    16 = @(993.55,    978.45)
    17 = @(631,       633.75)
    18 = @(2302.3,    2269.9)
    19 = @(266.9,     263.75)
    20 = @(22287.05,  22308.7)
    21 = @(361.35,    355.8)
    22 = @(820.3,     818.2)
    23 = @(660.85,    668.7)
    24 = @(947.3,     964.65)
    25 = @(431.3,     430.35)
    26 = @(165.6,     164.95)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
}

# Row 24 (TM) now uses the same number format / style as the other
# highlighted rows (e.g. B3) instead of the default style.
$ws.Range("B24").NumberFormat = $ws.Range("B3").NumberFormat

# Move the saved selection from K11:L12 to I12.
$ws.Range("I12").Select() | Out-Null
